$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:D1) text - add "Graph " prefix ---
$ws.Range("A1").Value2 = "Graph BFS Memory Usage 10000 (bytes)"
$ws.Range("B1").Value2 = "Graph DFS Memory Usage 10000 (bytes)"
$ws.Range("C1").Value2 = "Graph BFS Memory Usage 1000 (bytes)"
$ws.Range("D1").Value2 = "Graph DFS Memory Usage 1000 (bytes)"

# --- Add new header columns E1:H1 for Tree data ---
$ws.Range("E1").Value2 = "Tree BFS Memory Usage 10000 (bytes)"
$ws.Range("F1").Value2 = "Tree DFS Memory Usage 10000 (bytes)"
$ws.Range("G1").Value2 = "Tree BFS Memory Usage 1000 (bytes)"
$ws.Range("H1").Value2 = "Tree DFS Memory Usage 1000 (bytes)"

# --- Remove old rows 3-13 (only rows 1 and 2 remain from before) ---
$ws.Range("A3:D13").EntireRow.Delete()

# --- Update row 2 values (Graph columns) and add Tree data (E2:H2) ---
$ws.Range("A2").Value2 = 9312680
$ws.Range("B2").Value2 = 7982304
$ws.Range("C2").Value2 = 1316056
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = 6658600
$ws.Range("F2").Value2 = 5326808
$ws.Range("G2").Value2 = 2663432
$ws.Range("H2").Value2 = 0

# --- New row 3: only G3 and H3 populated ---
$ws.Range("G3").Value2 = 0
$ws.Range("H3").Value2 = 1316048

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 33.166666666666664
$ws.Columns.Item(2).ColumnWidth = 33.333333333333336
$ws.Columns.Item(3).ColumnWidth = 32.5
$ws.Columns.Item(4).ColumnWidth = 32.666666666666664
$ws.Columns.Item(5).ColumnWidth = 32.166666666666664
$ws.Columns.Item(6).ColumnWidth = 32.166666666666664
$ws.Columns.Item(7).ColumnWidth = 31.0
$ws.Columns.Item(8).ColumnWidth = 31.333333333333332

# --- Reset selection to A1 (closest approximation to default view state) ---
[void]$ws.Range("A1").Select()

Write-Host "All edits applied"
